$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1033
$ws.Range("I53").Value = 866.25
$ws.Range("J53").Value = 1700
$ws.Range("K53").Value = 866.25
$ws.Range("L53").Value = 1700
$ws.Range("M53").Value = -229.25
$ws.Range("N53").Value = -2974

$ws.Range("H76").Value = 3979.2778
$ws.Range("J76").Value = 7316
$ws.Range("L76").Value = 7316
$ws.Range("N76").Value = -7946

$ws.Range("H79").Value = 3979.2778
$ws.Range("J79").Value = 7316
$ws.Range("L79").Value = 7316
$ws.Range("N79").Value = -9500

$ws.Range("H86").Value = 111113000
$ws.Range("I86").Value = 250001550
$ws.Range("K86").Value = 250001550
$ws.Range("M86").Value = -250000427

$ws.Range("H87").Value = 22222
$ws.Range("J87").Value = 22222
$ws.Range("L87").Value = 22222
$ws.Range("N87").Value = -24718

$ws.Range("H88").Value = 2169.875
$ws.Range("I88").Value = 1374.5714
$ws.Range("J88").Value = 2788.4443
$ws.Range("K88").Value = 1374.5714
$ws.Range("L88").Value = 2788.4443
$ws.Range("M88").Value = -968.5714
$ws.Range("N88").Value = -3600.4443

$ws.Range("H89").Value = 111113000
$ws.Range("I89").Value = 250001550
$ws.Range("K89").Value = 1250007750
$ws.Range("M89").Value = -1250002134

$ws.Range("H90").Value = 22222
$ws.Range("J90").Value = 22222
$ws.Range("L90").Value = 66666
$ws.Range("N90").Value = -79146

$ws.Range("H91").Value = 2169.875
$ws.Range("I91").Value = 1374.5714
$ws.Range("J91").Value = 2788.4443
$ws.Range("K91").Value = 1374.5714
$ws.Range("L91").Value = 2788.4443
$ws.Range("M91").Value = 29.42859999999996
$ws.Range("N91").Value = -5596.4443

$ws.Range("H116").Value = 11015.519
$ws.Range("I116").Value = 11809.4375
$ws.Range("J116").Value = 9860.727999999999
$ws.Range("K116").Value = 11809.4375
$ws.Range("L116").Value = 9860.727999999999
$ws.Range("M116").Value = -8367.4375
$ws.Range("N116").Value = -16744.728

$ws.Range("H125").Value = 5535.1
$ws.Range("I125").Value = 4992.6665
$ws.Range("J125").Value = 6348.75
$ws.Range("K125").Value = 44933.9985
$ws.Range("L125").Value = 57138.75
$ws.Range("M125").Value = -42473.9985
$ws.Range("N125").Value = -62058.75

$ws.Range("H135").Value = 11399.4
$ws.Range("I135").Value = 7332.3335
$ws.Range("K135").Value = 65991.0015
$ws.Range("M135").Value = -63456.0015

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1275.0869
$ws.Range("I32").Value = 551.59015
$ws.Range("J32").Value = 6791.75
$ws.Range("K32").Value = 551.59015
$ws.Range("L32").Value = 6791.75
$ws.Range("M32").Value = -264.59015
$ws.Range("N32").Value = -7365.75

$ws.Range("H61").Value = 2879.1
$ws.Range("I61").Value = 2175.6875
$ws.Range("J61").Value = 3683
$ws.Range("K61").Value = 2175.6875
$ws.Range("L61").Value = 3683
$ws.Range("M61").Value = -1963.6875
$ws.Range("N61").Value = -4107

$ws.Range("H74").Value = 2513.4443
$ws.Range("I74").Value = 1988.7142
$ws.Range("J74").Value = 4350
$ws.Range("K74").Value = 1988.7142
$ws.Range("L74").Value = 4350
$ws.Range("M74").Value = -1114.7142
$ws.Range("N74").Value = -6098

$ws.Range("H77").Value = 2513.4443
$ws.Range("I77").Value = 1988.7142
$ws.Range("J77").Value = 4350
$ws.Range("K77").Value = 9943.571
$ws.Range("L77").Value = 21750
$ws.Range("M77").Value = -5575.571
$ws.Range("N77").Value = -30486

$ws.Range("H102").Value = 5067.2593
$ws.Range("J102").Value = 14996.333
$ws.Range("L102").Value = 14996.333
$ws.Range("N102").Value = -18240.333

$ws.Range("H122").Value = 2384.0386
$ws.Range("I122").Value = 1548.7059
$ws.Range("K122").Value = 4646.1177
$ws.Range("M122").Value = -2196.1177

$ws.Range("H136").Value = 2879.1
$ws.Range("I136").Value = 2175.6875
$ws.Range("J136").Value = 3683
$ws.Range("K136").Value = 6527.0625
$ws.Range("L136").Value = 11049
$ws.Range("M136").Value = -3977.0625
$ws.Range("N136").Value = -16149

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 49024.5
$ws.Range("J110").Value = 49024.5
$ws.Range("L110").Value = 49024.5
$ws.Range("N110").Value = -57204.5

$ws.Range("H134").Value = 11053.421
$ws.Range("I134").Value = 4713.846
$ws.Range("J134").Value = 24789.166
$ws.Range("K134").Value = 14141.538
$ws.Range("L134").Value = 74367.49800000001
$ws.Range("M134").Value = -11606.538
$ws.Range("N134").Value = -79437.49800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws.Range("H31").Value = 2683.3
$ws.Range("I31").Value = 1399.3334
$ws.Range("J31").Value = 2909.8823
$ws.Range("K31").Value = 1399.3334
$ws.Range("L31").Value = 2909.8823
$ws.Range("M31").Value = -1104.3334
$ws.Range("N31").Value = -3499.8823

$ws.Range("H34").Value = 2683.3
$ws.Range("I34").Value = 1399.3334
$ws.Range("J34").Value = 2909.8823
$ws.Range("K34").Value = 1399.3334
$ws.Range("L34").Value = 2909.8823
$ws.Range("M34").Value = -1197.3334
$ws.Range("N34").Value = -3313.8823

$ws.Range("H86").Value = 10801.52
$ws.Range("I86").Value = 10328.071
$ws.Range("J86").Value = 11404.091
$ws.Range("K86").Value = 10328.071
$ws.Range("L86").Value = 11404.091
$ws.Range("M86").Value = -9205.071
$ws.Range("N86").Value = -13650.091

$ws.Range("H89").Value = 10801.52
$ws.Range("I89").Value = 10328.071
$ws.Range("J89").Value = 11404.091
$ws.Range("K89").Value = 51640.355
$ws.Range("L89").Value = 57020.455
$ws.Range("M89").Value = -46024.355
$ws.Range("N89").Value = -68252.455

$ws.Range("H132").Value = 4865.2666
$ws.Range("I132").Value = 3564.8696
$ws.Range("K132").Value = 10694.6088
$ws.Range("M132").Value = -8164.6088

$ws.Range("H134").Value = 9181.299999999999
$ws.Range("I134").Value = 3028.4285
$ws.Range("J134").Value = 11053.913
$ws.Range("K134").Value = 9085.2855
$ws.Range("L134").Value = 33161.739
$ws.Range("M134").Value = -6550.2855
$ws.Range("N134").Value = -38231.739

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 296860.28
$ws.Range("J68").Value = 359709
$ws.Range("L68").Value = 1079127
$ws.Range("N68").Value = -1080749

$ws.Range("H71").Value = 296860.28
$ws.Range("J71").Value = 359709
$ws.Range("L71").Value = 3237381
$ws.Range("N71").Value = -3245493

$ws.Range("H132").Value = 4172.125
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 4172.125
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 37549.125
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -42609.125

$ws.Range("H134").Value = 6230.143
$ws.Range("I134").Value = 6230.143
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 18690.429
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -13620.429
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6322.4
$ws.Range("I102").Value = 5211.8335
$ws.Range("J102").Value = 7988.25
$ws.Range("K102").Value = 5211.8335
$ws.Range("L102").Value = 7988.25
$ws.Range("M102").Value = -3589.8335
$ws.Range("N102").Value = -11232.25

$ws.Range("H122").Value = 5378.32
$ws.Range("I122").Value = 2373.25
$ws.Range("J122").Value = 10720.667
$ws.Range("K122").Value = 7119.75
$ws.Range("L122").Value = 32162.001
$ws.Range("M122").Value = -4669.75
$ws.Range("N122").Value = -37062.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4129.75
$ws.Range("I82").Value = 1980.7
$ws.Range("J82").Value = 14875
$ws.Range("K82").Value = 1980.7
$ws.Range("L82").Value = 14875
$ws.Range("M82").Value = -1619.7
$ws.Range("N82").Value = -15597

$ws.Range("H85").Value = 4129.75
$ws.Range("I85").Value = 1980.7
$ws.Range("J85").Value = 14875
$ws.Range("K85").Value = 1980.7
$ws.Range("L85").Value = 14875
$ws.Range("M85").Value = -732.7
$ws.Range("N85").Value = -17371

$ws.Range("H93").Value = 7060.4165
$ws.Range("I93").Value = 6172.1
$ws.Range("J93").Value = 11502
$ws.Range("K93").Value = 6172.1
$ws.Range("L93").Value = 11502
$ws.Range("M93").Value = -4924.1
$ws.Range("N93").Value = -13998

$ws.Range("H132").Value = 2542.05
$ws.Range("J132").Value = 3317.1667
$ws.Range("L132").Value = 9951.500100000001
$ws.Range("N132").Value = -15011.5001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H96").Value = 2838.3333
$ws.Range("I96").Value = 2099.5
$ws.Range("J96").Value = 2952
$ws.Range("K96").Value = 2099.5
$ws.Range("L96").Value = 2952
$ws.Range("M96").Value = -726.5
$ws.Range("N96").Value = -5698

$ws.Range("H100").Value = 1188
$ws.Range("I100").Value = 907.4545000000001
$ws.Range("K100").Value = 1814.909
$ws.Range("M100").Value = -1273.909
